$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update 想去人数 (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 197
$ws1.Range("F4").Value = 2304
$ws1.Range("F5").Value = 1735
$ws1.Range("F6").Value = 332
$ws1.Range("F8").Value = 793
$ws1.Range("F9").Value = 163

# Sheet "全部类型" (sheet4): update 想去人数 (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 197
$ws4.Range("F4").Value = 2304
$ws4.Range("F5").Value = 1735
$ws4.Range("F6").Value = 332
$ws4.Range("F9").Value = 793
$ws4.Range("F10").Value = 163
